$wb = $excel.ActiveWorkbook

# Map of F-column cell updates (row -> new value) that apply identically
# to the "展览" and "全部类型" worksheets.
$updates = @{
    "F5"  = 2563
    "F9"  = 1309
    "F13" = 1158
    "F15" = 323
    "F21" = 2346
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
